$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new transaction was recorded on 2024-09-05 16:22:23 ("transfer").
# This pushes every existing September/August row down by one, so insert
# a fresh row above the current row 31 (shifts 31..70 -> 32..71, and
# grows the used range to A1:Y71), then populate the new row's
# September_Details / September_Date cells.
$ws.Rows.Item(31).Insert()

$ws.Range("R31").Value = "transfer"
$ws.Range("S31").Value = "2024-09-05 16:22:23"
